# Apply unit conversion (Sublimation Energy, row 4) and
# number-of-molecule correction (lattice parameters a/b/c, rows 5-7)
# to both worksheets ("N-Ac,N'Me-Ala" and "Adipamide").

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Row 4 - Sublimation Energy: unit conversion (divide by 4.184, cal -> kJ)
    $ws.Range("B4").Value = 26.2297155823325
    $ws.Range("D4").Value = 26.2297155823325
    $ws.Range("F4").Value = 26.2297155823325

    # Row 5 - a: number of molecule correction (x10)
    $ws.Range("B5").Value = 9.600429514221128
    $ws.Range("D5").Value = 9.600429514221128
    $ws.Range("F5").Value = 9.600429514221128

    # Row 6 - b: number of molecule correction (x10)
    $ws.Range("B6").Value = 6.713888386566135
    $ws.Range("D6").Value = 6.713888386566135
    $ws.Range("F6").Value = 6.713888386566135
    $ws.Range("C6").Value = 0.4990181765615594
    $ws.Range("E6").Value = 0.4990181765615594
    $ws.Range("G6").Value = 0.4990181765615594

    # Row 7 - c: number of molecule correction (x10)
    $ws.Range("B7").Value = 7.260101054623197
    $ws.Range("D7").Value = 7.260101054623197
    $ws.Range("F7").Value = 7.260101054623197
    $ws.Range("C7").Value = 0.3245547642408042
    $ws.Range("E7").Value = 0.3245547642408042
    $ws.Range("G7").Value = 0.3245547642408042
}
